# update the demo templates a bit more
$wb = $excel.ActiveWorkbook

# --- site sheet ---
$site = $wb.Worksheets.Item("site")
$site.Activate()
$site.Range("B6").Value = $true
$site.Range("B7").Select()

# --- outing sheet ---
$outing = $wb.Worksheets.Item("outing")
$outing.Activate()
$outing.Range("B6").Value = $true
$outing.Range("C10").Select()

# --- recapture sheet ---
$recapture = $wb.Worksheets.Item("recapture")
$recapture.Activate()
$recapture.Range("E8").Select()

# --- capture sheet (ends as the active tab, per workbook activeTab=2) ---
$capture = $wb.Worksheets.Item("capture")
$capture.Activate()
$capture.Range("L6").Value = $true
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$capture.Range("F14").Select()
